$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.475.69"
$ws.Range("E2").Value = "  +2.78%  "
$ws.Range("D3").Value = "1.730.18"
$ws.Range("E3").Value = "  +3.18%  "
$ws.Range("D4").Value = "'0.9993"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'244.00"
$ws.Range("E5").Value = "  +2.86%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'0.4800"
$ws.Range("E7").Value = "  +3.74%  "
$ws.Range("D8").Value = "'0.2671"
$ws.Range("E8").Value = "  +2.97%  "
$ws.Range("D9").Value = "'0.06229"
$ws.Range("E9").Value = "  +1.46%  "
$ws.Range("D10").Value = "1.732.20"
$ws.Range("E10").Value = "  +3.60%  "
$ws.Range("D11").Value = "'0.07113"
$ws.Range("E11").Value = "  +1.64%  "
$ws.Range("E12").Value = "  +5.52%  "
$ws.Range("D13").Value = "'0.6180"
$ws.Range("E13").Value = "  +7.22%  "
$ws.Range("E14").Value = "  +4.20%  "
$ws.Range("D15").Value = "'77.07"
$ws.Range("E15").Value = "  +2.26%  "
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").Value = "26.486.38"
$ws.Range("E17").Value = "  +2.84%  "
$ws.Range("D18").Value = "'0.9996"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").Value = "'0.000006928"
$ws.Range("E19").Value = "  +3.44%  "
$ws.Range("E20").Value = "  +2.67%  "
$ws.Range("D21").Value = "1.954.09"
$ws.Range("E21").Value = "  +3.62%  "
$ws.Range("D22").Value = "'4.563"
$ws.Range("E22").Value = "  +2.23%  "
$ws.Range("D23").Value = "'8.908"
$ws.Range("E23").Value = "  +2.82%  "
$ws.Range("D24").Value = "'5.316"
$ws.Range("E24").Value = "  +1.60%  "
$ws.Range("D25").Value = "'136.40"
$ws.Range("E25").Value = "  +1.62%  "
$ws.Range("D26").Value = "'15.34"
$ws.Range("E26").Value = "  +2.33%  "
$ws.Range("E27").Value = "  +4.32%  "
$ws.Range("D28").Value = "'1.412"
$ws.Range("E28").Value = "  +1.58%  "
$ws.Range("D29").Value = "'106.59"
$ws.Range("E29").Value = "  +1.91%  "
$ws.Range("D30").Value = "'3.989"
$ws.Range("E30").Value = "  +1.11%  "
$ws.Range("D31").Value = "'0.07986"
$ws.Range("E31").Value = "  +4.02%  "
$ws.Range("D32").Value = "'3.745"
$ws.Range("E32").Value = "  +3.60%  "
$ws.Range("D33").Value = "'0.04563"
$ws.Range("E33").Value = "  +5.24%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'2.613"
$ws.Range("E34").Value = "  +0.44%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.6409"
$ws.Range("E35").Value = "  +5.30%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'0.9907"
$ws.Range("E36").Value = "  +4.13%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "'0.9460"
$ws.Range("E37").Value = "  +1.22%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'1.988"
$ws.Range("E38").Value = "  +6.95%  "
$ws.Range("B39").Value = "Quant"
$ws.Range("C39").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D39").Value = "'107.53"
$ws.Range("E39").Value = "  -1.59%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.407"
$ws.Range("E40").Value = "  -1.48%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "'1.006"
$ws.Range("E41").Value = "  +0.76%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.01503"
$ws.Range("E42").Value = "  +3.75%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.652"
$ws.Range("E43").Value = "  +11.72%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.3910"
$ws.Range("E44").Value = "  +4.97%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "'6.934"
$ws.Range("E45").Value = "  +12.81%  "
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "'0.1192"
$ws.Range("E46").Value = "  +6.68%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.05330"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'7.897"
$ws.Range("E48").Value = "  +3.43%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'30.81"
$ws.Range("E49").Value = "  -2.17%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'1.271"
$ws.Range("E50").Value = "  +5.01%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "'0.3428"
$ws.Range("E51").Value = "  +3.28%  "
